$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: update a cell's text while keeping it stored as Text.
# Many of the source values look numeric (e.g. "214.92"), so Excel would
# otherwise silently convert them to real numbers on assignment. We force
# the Text number format right before assigning, then restore the default
# "Normal" style so no stray formatting is left behind on the cell.
function Set-CellText {
    param($cellRef, $text)
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-CellText "D2" "25.846.30"
Set-CellText "E2" "  +0.11%  "
Set-CellText "D3" "1.634.18"
Set-CellText "E3" "  +0.31%  "
Set-CellText "E4" "  -0.12%  "
Set-CellText "D5" "214.92"
Set-CellText "E5" "  -0.28%  "
Set-CellText "D6" "0.504"
Set-CellText "E6" "  -0.47%  "
Set-CellText "E7" "  -0.14%  "
Set-CellText "D8" "0.257"
Set-CellText "E8" "  -0.26%  "
Set-CellText "E9" "  -0.15%  "
Set-CellText "D10" "19.87"
Set-CellText "E10" "  +2.70%  "
Set-CellText "E11" "  +0.04%  "
Set-CellText "E12" "  -0.06%  "
Set-CellText "D13" "1.860.72"
Set-CellText "E13" "  +0.37%  "
Set-CellText "D14" "1.631.36"
Set-CellText "E14" "  +0.49%  "
Set-CellText "E15" "  +0.32%  "
Set-CellText "D16" "0.0₃0766"
Set-CellText "E16" "  +1.88%  "
Set-CellText "D17" "63.02"
Set-CellText "E17" "  -0.40%  "
Set-CellText "D18" "25.859.34"
Set-CellText "D20" "193.43"
Set-CellText "E20" "  +0.38%  "
Set-CellText "D21" "4.39"
Set-CellText "E21" "  +2.13%  "
Set-CellText "D23" "6.19"
Set-CellText "E23" "  +3.11%  "
Set-CellText "E24" "  -0.18%  "
Set-CellText "D25" "1.76"
Set-CellText "E25" "  -4.22%  "
Set-CellText "D26" "139.24"
Set-CellText "E26" "  -0.75%  "
Set-CellText "E27" "  -4.40%  "
Set-CellText "D28" "6.81"
Set-CellText "E28" "  +1.43%  "
Set-CellText "E29" "  +0.76%  "
Set-CellText "E30" "  -0.02%  "
Set-CellText "E31" "  +1.86%  "
Set-CellText "E33" "  +1.81%  "
Set-CellText "D34" "1.57"
Set-CellText "E34" "  +0.87%  "
Set-CellText "D35" "2.38"
Set-CellText "E35" "  +0.13%  "
Set-CellText "E36" "  +0.76%  "
Set-CellText "E37" "  +1.11%  "
Set-CellText "D38" "1.122.52"
Set-CellText "E38" "  -1.22%  "
Set-CellText "E39" "  +0.57%  "
Set-CellText "D40" "0.0156"
Set-CellText "E40" "  +0.09%  "
Set-CellText "E41" "  +0.09%  "
Set-CellText "E42" "  -1.08%  "
Set-CellText "D43" "99.62"
Set-CellText "E44" "  +0.68%  "
Set-CellText "E45" "  -2.51%  "
Set-CellText "D46" "55.41"
Set-CellText "E46" "  +0.85%  "
Set-CellText "E47" "  -4.96%  "
Set-CellText "D49" "7.64"
Set-CellText "E49" "  +0.77%  "
Set-CellText "E50" "  +0.14%  "
Set-CellText "E51" "  +6.63%  "

